$wb = $excel.ActiveWorkbook

# --- Sheet "Overview": row 5 (846b33cd-...md) status changes from
#     "Ready for handoff" to "Handed back: in sync with en-US"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B5").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C5").Value = "Handed back: in sync with en-US"

# --- Sheet "zh-cn": same status change, plus newly-populated handback
#     columns (Latest Target File / Latest Handback File / Latest
#     Handback DateTime) for row 5, now that the file has been handed
#     back.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B5").Value = "Handed back: in sync with en-US"
$wsZh.Range("E5").Value = "846b33cd-68a3-4ac1-8cff-73eae6dd70ea.md"
$wsZh.Range("F5").Value = "846b33cd-68a3-4ac1-8cff-73eae6dd70ea.c6ffe1ff193545ef94f852a13dc9edcf42ca0990.zh-cn.xlf"
$wsZh.Range("G5").Value = "2016-02-22 17:45:37"

$wsZh.Hyperlinks.Add(
    $wsZh.Range("E5"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/fa53bbda021166b2e977ed529487d43c9b32cbd6/e2e/846b33cd-68a3-4ac1-8cff-73eae6dd70ea.md",
    "",
    "",
    "846b33cd-68a3-4ac1-8cff-73eae6dd70ea.md"
) | Out-Null

$wsZh.Hyperlinks.Add(
    $wsZh.Range("F5"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/18029b537c0d7e9dd221b93c704c4d1f63cf0b27/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/846b33cd-68a3-4ac1-8cff-73eae6dd70ea.c6ffe1ff193545ef94f852a13dc9edcf42ca0990.zh-cn.xlf",
    "",
    "",
    "846b33cd-68a3-4ac1-8cff-73eae6dd70ea.c6ffe1ff193545ef94f852a13dc9edcf42ca0990.zh-cn.xlf"
) | Out-Null

# --- Sheet "de-de": same treatment
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B5").Value = "Handed back: in sync with en-US"
$wsDe.Range("E5").Value = "846b33cd-68a3-4ac1-8cff-73eae6dd70ea.md"
$wsDe.Range("F5").Value = "846b33cd-68a3-4ac1-8cff-73eae6dd70ea.c6ffe1ff193545ef94f852a13dc9edcf42ca0990.de-de.xlf"
$wsDe.Range("G5").Value = "2016-02-22 17:45:58"

$wsDe.Hyperlinks.Add(
    $wsDe.Range("E5"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/4a873735bd514d53eff3efc79fae566d313a7654/e2e/846b33cd-68a3-4ac1-8cff-73eae6dd70ea.md",
    "",
    "",
    "846b33cd-68a3-4ac1-8cff-73eae6dd70ea.md"
) | Out-Null

$wsDe.Hyperlinks.Add(
    $wsDe.Range("F5"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ee5da840b5359054cef11cb32f56a459cd08a7ed/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/846b33cd-68a3-4ac1-8cff-73eae6dd70ea.c6ffe1ff193545ef94f852a13dc9edcf42ca0990.de-de.xlf",
    "",
    "",
    "846b33cd-68a3-4ac1-8cff-73eae6dd70ea.c6ffe1ff193545ef94f852a13dc9edcf42ca0990.de-de.xlf"
) | Out-Null
